$d = $word.ActiveDocument

# The second paragraph holds a single merge field " m:self.name " (with
# "self" highlighted in accent6/darker-25%). Replace the field with plain
# literal runs "{", "m", ":", "self" (still highlighted) and ".name}" so
# the field codes become ordinary template-engine text instead of a real
# Word field.

$field = $d.Fields.Item(1)
$para = $d.Paragraphs.Item(2)
$insertPos = $para.Range.Start

# Remove the whole field (fldChar begin/end + all instrText runs).
$field.Delete()

$target = $d.Range($insertPos, $insertPos)

$runsXml = '<w:r><w:t>{</w:t></w:r>' + `
           '<w:r><w:t>m</w:t></w:r>' + `
           '<w:r><w:t>:</w:t></w:r>' + `
           '<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>self</w:t></w:r>' + `
           '<w:r><w:t xml:space="preserve">.name}</w:t></w:r>'

$packageXml = '<?xml version="1.0" standalone="yes"?>' + `
              '<?mso-application progid="Word.Document"?>' + `
              '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
              '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
              '<pkg:xmlData>' + `
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
              '<w:body><w:p>' + $runsXml + '</w:p></w:body>' + `
              '</w:document>' + `
              '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($packageXml)
